# Change copyright notice from "Inc." to "PBC" in the footer text boxes
# (e.g. "RStudio, Inc.  *  " -> "RStudio, PBC  *  ") on every slide of the
# presentation, wherever it occurs, while leaving everything else in the
# text (including the bullet separator that follows) untouched.

$p = $ppt.ActivePresentation

$oldText = ", Inc."
$newText = ", PBC"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        try {
            if ($shape.HasTextFrame -eq $false) { continue }

            $tf = $shape.TextFrame
            if ($tf.HasText -eq $false) { continue }

            $tr = $tf.TextRange
            if ($tr.Text -eq $null -or $tr.Text.Length -eq 0) { continue }

            # Keep searching/replacing in case the same shape contains the
            # phrase more than once.
            $searchStart = 0
            while ($true) {
                $found = $tr.Find($oldText, $searchStart)
                if ($found -eq $null) { break }

                # Select exactly the characters that make up ", Inc." and
                # replace only that portion with ", PBC". This preserves the
                # run's existing character formatting (font, size, color,
                # bold/strike, etc.) as well as everything else in the text
                # box (e.g. the trailing "  " + bullet + "  " that follows).
                $target = $tr.Characters($found.Start, $found.Length)
                $target.Text = $newText

                $searchStart = $found.Start + $newText.Length - 1
                if ($searchStart -lt 0) { $searchStart = 0 }
            }
        } catch {
            # Ignore shapes that don't support text access and keep going.
            continue
        }
    }
}
